# "player fire moved into playerclass"
# The PLAYER_FIRE related UI-parameter blocks (button3 / button4, rows 37-48)
# are removed from the player_parameter sheet - their cell content is cleared
# out (styles/formatting are left intact), since that configuration now lives
# in the player class instead of this datasheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("player_parameter")

# Clear the two now-obsolete parameter blocks (rows 37-42 = "button3",
# rows 43-48 = "button4") while keeping the pre-existing cell styles.
$ws.Range("A37:B48").ClearContents()

# Reflect the final on-screen selection/scroll position left behind by the edit.
$ws.Activate()
$ws.Range("I43").Select()
$excel.ActiveWindow.ScrollRow = 31
$excel.ActiveWindow.ScrollColumn = 1
